# Assignment-1 (Udanous_) (Recovered).xlsx -- apply the captured edit
#
# Summary of changes (see xml_diff):
#  1. Workbook view: active tab moves from "Q1 to Q10" (index 1) to "Q19" (index 3);
#     window is no longer minimized (that flag disappears on any re-save of this engine).
#  2. "Q1 to Q10" sheet: no longer the tab-selected sheet; its cursor selection moves
#     from H275 to F45.
#  3. "Q1 to Q10" formula edits:
#       E16  COUNTIF(C3:C12,">100")  -> COUNTIF(C3:C12,">=100")
#       B26  MONTH(B23)              -> DATEDIF(B23,TODAY(),"M")
#       D54  COUNT(D45:D53)          -> ROWS(D45:D53)
#     (B27/C27/B81/B341 keep their formulas; their cached values merely reflect a later
#      TODAY() recalculation, handled automatically by the post-script recalc.)
#  4. "Q19" sheet becomes the tab-selected sheet; cursor selection moves to J23.
#  5. "Q19" sheet: columns C:J (rows 2-20) get the age-calculation formulas that were
#     previously blank placeholder cells.

$wb = $excel.ActiveWorkbook

$wsQ1 = $wb.Worksheets.Item("Q1 to Q10")
$wsQ19 = $wb.Worksheets.Item("Q19")

# ---------------------------------------------------------------------------
# "Q1 to Q10" formula fixes
# ---------------------------------------------------------------------------

# >100  ->  >=100
$wsQ1.Range("E16").Formula = "=COUNTIF(C3:C12,"">=100"")"

# MONTH(B23)  ->  DATEDIF(B23,TODAY(),"M")
$wsQ1.Range("B26").Formula = "=DATEDIF(B23,TODAY(),""M"")"

# COUNT(D45:D53)  ->  ROWS(D45:D53)
$wsQ1.Range("D54").Formula = "=ROWS(D45:D53)"

# ---------------------------------------------------------------------------
# "Q19" sheet: fill in the age-calculation table (previously empty C:J)
# ---------------------------------------------------------------------------

# Row 2 uses standalone (non-shared) formulas referencing row 2.
$wsQ19.Range("C2").Formula = "=DATEDIF(B2,TODAY(),""Y"")"
$wsQ19.Range("D2").Formula = "=DATEDIF(B2,TODAY(),""M"")"
$wsQ19.Range("E2").Formula = "=DATEDIF(B2,TODAY(),""D"")"
$wsQ19.Range("F2").Formula = "=DATEDIF(B2,TODAY(),""YM"")"
$wsQ19.Range("G2").Formula = "=DATEDIF(B2,TODAY(),""MD"")"
$wsQ19.Range("H2").Formula = "=DATEDIF(B2,TODAY(),""YD"")"
$wsQ19.Range("I2").Formula = "=CONCATENATE(C2,"" Years"","" "",G2,"" "",""Days"")"
$wsQ19.Range("J2").Formula = "=CONCATENATE(C2,"" Years"","" "",F2,"" Months"","" "",G2,"" Days"")"

# Rows 3-20 fill as shared formulas (one Formula assignment per column range, like an
# Excel fill-down, so the engine emits t="shared" groups the same way Excel does).
$wsQ19.Range("C3:C20").Formula = "=DATEDIF(B3,TODAY(),""Y"")"
$wsQ19.Range("D3:D20").Formula = "=DATEDIF(B3,TODAY(),""M"")"
$wsQ19.Range("E3:E20").Formula = "=DATEDIF(B3,TODAY(),""D"")"
$wsQ19.Range("F3:F20").Formula = "=DATEDIF(B3,TODAY(),""YM"")"
$wsQ19.Range("G3:G20").Formula = "=DATEDIF(B3,TODAY(),""MD"")"
$wsQ19.Range("H3:H20").Formula = "=DATEDIF(B3,TODAY(),""YD"")"
$wsQ19.Range("I3:I20").Formula = "=CONCATENATE(C3,"" Years"","" "",G3,"" "",""Days"")"
$wsQ19.Range("J3:J20").Formula = "=CONCATENATE(C3,"" Years"","" "",F3,"" Months"","" "",G3,"" Days"")"

# ---------------------------------------------------------------------------
# View / selection / active-tab changes
# ---------------------------------------------------------------------------

# Move the cursor on "Q1 to Q10" (this also de-selects it as the active tab once
# another sheet is activated below).
$wsQ1.Range("F45").Select()

# "Q19" becomes the active sheet/tab, with its cursor on J23.
$wsQ19.Activate()
$wsQ19.Range("J23").Select()
